# DGS risk-matrix time series: append the 2021/08/27 report row.
#
# The sheet holds a simple time series in columns A:E (date, two raw
# counts and two ratios). Each existing row lives in row 2..72; the new
# report adds one more row (row 73) with the same column layout/styles.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 73

# Column A carries the report date but is stored as *text* (all the
# other date cells in this column are shared-string text, not real
# date serials) even though the column uses a yyyy/mm/dd number format.
# Excel's normal text auto-recognizes "2021/08/27" as a date and would
# silently convert the cell to a date serial, so force a text format
# first, enter the literal value, then restore the column's usual
# yyyy/mm/dd display format (the stored value remains text either way).
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2021/08/27"
$dateCell.NumberFormat = "yyyy/mm/dd"

$ws.Cells.Item($newRow, 2).Value = 312.7
$ws.Cells.Item($newRow, 3).Value = 317.7
$ws.Cells.Item($newRow, 4).Value = 0.99
$ws.Cells.Item($newRow, 5).Value = 0.99

# Match the author's cursor position after typing the new row: one row
# below the freshly-entered data.
$ws.Range("A74").Select()
